$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the realization value for row 6 (KAP Budi) from the dashboard
$ws.Range("C6").Value = 300000000

# Add REMAINING formula (column D) = CONTRACT_VALUE - REALIZATION
$ws.Range("D2").Formula = "=B2-C2"
$ws.Range("D3:D24").Formula = "=B3-C3"

# Add REALIZED_PCT formula (column E) = (REALIZATION / CONTRACT_VALUE) * 100
$ws.Range("E2").Formula = "=(C2/B2)*100"
$ws.Range("E3:E24").Formula = "=(C3/B3)*100"

# Move the active selection to C10, matching the dashboard's last-touched cell
$ws.Range("C10").Select()
